# Fix data for the "Centre A Ben Mansour" group (rows 2-4):
# - Align the latitude/longitude of rows 3 and 4 with row 2 (they were
#   accidental near-duplicate coordinates with typo'd longitude values).
# - Update number_of_kiosks (column L) for every registration center row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Align latitude (G) for rows 3 and 4 with row 2
$ws.Range("G3").Value = 34.521169999999998
$ws.Range("G4").Value = 34.521169999999998

# Align longitude (H) for rows 3 and 4 with row 2's value. The longitude is
# stored as text (it starts with a non-breaking space), so copy/paste the
# value from H2 instead of retyping it, which keeps the text formatting
# intact without touching the cell style.
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H3").PasteSpecial(-4163) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# number_of_kiosks (column L): rows 2-4 -> 3, rows 5-46 -> 2
$ws.Range("L2").Value = 3
$ws.Range("L3").Value = 3
$ws.Range("L4").Value = 3

for ($r = 5; $r -le 46; $r++) {
    $ws.Cells.Item($r, 12).Value = 2
}

# Reflect the last interactive selection recorded in the workbook
$ws.Range("L28").Select()
